# Insert a new data row above current row 269 ("Vega Modelo de Temuco" /
# Albahaca weekly price sheet). All existing rows 269:321 shift down to
# 270:322, and the new row 269 holds a fresh weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 269 (and everything below it) down by one row.
$ws.Rows.Item(269).Insert()

# Populate the newly inserted row 269 with the new record.
$ws.Cells.Item(269, 1).Value  = 10
$ws.Cells.Item(269, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(269, 3).Value  = "La Araucanía"
$ws.Cells.Item(269, 4).Value  = 44951
$ws.Cells.Item(269, 5).Value  = 9
$ws.Cells.Item(269, 6).Value  = 100112052
$ws.Cells.Item(269, 7).Value  = "Albahaca"
$ws.Cells.Item(269, 8).Value  = "Sin especificar"
$ws.Cells.Item(269, 9).Value  = "Primera"
$ws.Cells.Item(269, 10).Value = 25
$ws.Cells.Item(269, 11).Value = 6000
$ws.Cells.Item(269, 12).Value = 6000
$ws.Cells.Item(269, 13).Value = 6000
$ws.Cells.Item(269, 14).Value = "$/paquete"
$ws.Cells.Item(269, 15).Value = "Región del Maule"
$ws.Cells.Item(269, 16).Value = 6000
$ws.Cells.Item(269, 17).Value = 1
$ws.Cells.Item(269, 18).Value = "Hortaliza"

# Match the date formatting used by column D elsewhere on the sheet.
$ws.Cells.Item(269, 4).NumberFormat = $ws.Cells.Item(270, 4).NumberFormat
